# Edit described by the diff:
#  1. Apply strikethrough formatting to the first To-Do bullet ("Make sure
#     Goomba creature changes direction when hitting wall or bumping into
#     its friends.") to mark it as completed.
#  2. Relocate the auto-managed "_GoBack" bookmark from the end of the
#     "Implement turtle knocking over..." bullet to the end of the
#     "Create player killed animation..." bullet. Adding a new bookmark
#     with the same name removes the old one automatically, matching
#     Word's single-bookmark-per-name behavior.

$d = $word.ActiveDocument

# Find the paragraphs we need by their (stable) text content instead of a
# hard-coded index, so the script is robust to any unrelated reflow.
$goombaPara = $null
$killedAnimPara = $null
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "Make sure Goomba creature changes direction*" -and $goombaPara -eq $null) {
        $goombaPara = $p
    }
    if ($t -like "Create player killed animation when bumps into*" -and $killedAnimPara -eq $null) {
        $killedAnimPara = $p
    }
}

# --- 1. Strike through the first bullet item ---
$goombaPara.Range.Font.StrikeThrough = 1

# --- 2. Move the "_GoBack" bookmark to the end of $killedAnimPara ---
$insertionPoint = $killedAnimPara.Range.Duplicate
$insertionPoint.MoveEnd(1, -1) | Out-Null   # exclude the paragraph mark
$insertionPoint.Collapse(0)                 # collapse to its end (wdCollapseEnd = 0)

# A zero-length bookmark placed exactly at a paragraph-end boundary isn't
# handled correctly by this host, so temporarily insert a marker
# character, wrap the bookmark around it, then delete the marker. The
# bookmark collapses back to the same (now safe) position, landing right
# after the last run and before the paragraph mark -- exactly where the
# bookmark used to sit relative to its original (now previous) paragraph.
$insertionPoint.InsertAfter("|")
$d.Bookmarks.Add("_GoBack", $insertionPoint)

$markerRange = $insertionPoint.Duplicate
$markerRange.Collapse(1)                    # wdCollapseStart = 1
$markerRange.MoveEnd(1, 1) | Out-Null
$markerRange.Delete()

Write-Output "Done"
